$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Species Set")

# Delete the 31 "plant family" host-plant columns (BP:CT) which are no
# longer needed now that crops have been folded into the existing
# habitat categories. Everything to the right (old CU/CV -> new BP/BQ)
# shifts left to fill the gap.
$ws.Range("BP1:CT1").EntireColumn.Delete()

$ws.Range("A1:BQ28").AutoFilter() | Out-Null
